$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List of Projects")

# --- New rows: FizzBuzz (row 8) and Count words in a string (row 9) ---

# Row 8 - FizzBuzz
$ws.Range("B8").Value = "FizzBuzz"
$ws.Range("C8").Value = "prints the numbers from 1 to 100 for multiples of 3 writes fizz for multiples of 5 writes buzz for multiples of 3 and 5 writes fizzbuzz"
$ws.Range("D8").Value = "Text"
$ws.Range("E8").Value = 43519
$ws.Range("F8").Value = 43519
$ws.Range("G8").Value = "Java"

# Row 9 - Count words in a string (started, not finished)
$ws.Range("C9").Value = " Counts the number of individual words in a string. For added complexity read these strings in from a text file and generate a summary."
$ws.Range("B9").Value = "Count words in a string"
$ws.Range("D9").Value = "Text"
$ws.Range("E9").Value = 43519
$ws.Range("G9").Value = "Java"

# --- Date number format for new date cells (copy format from an existing
#     date cell so we reuse the existing built-in numFmt rather than create
#     a duplicate custom one) ---
$ws.Range("E4").Copy()
$ws.Range("E8:F8").PasteSpecial(-4122)
$ws.Range("E9").PasteSpecial(-4122)

# --- Column B: widen + wrap text (matches column C treatment) ---
$ws.Columns("B").ColumnWidth = 19.65
$ws.Columns("E").ColumnWidth = 14.65
$ws.Range("C4").Copy()
$ws.Range("B4:B9").PasteSpecial(-4122)

# Row heights for the two new rows (auto-fit equivalent, matches existing
# wrapped rows 5/6/7 which are also 43.5 / 29)
$ws.Rows(8).RowHeight = 43.5
$ws.Rows(9).RowHeight = 43.5

# --- Header tweaks: wrap the Project header, and add a date number format
#     to the Date Started / Date Finished headers ---
$ws.Range("B3").WrapText = $true
$ws.Range("E3:F3").NumberFormat = "m/d/yyyy"

# --- C9 description pasted from an external (GitHub) source: keeps its
#     own font color and wrap, but not the column's centered alignment ---
$ws.Range("C9").Font.Color = 2402084
$ws.Range("C9").VerticalAlignment = -4107
$ws.Range("C9").WrapText = $true

Write-Output "values set"
